$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SamplesTab query (B3): drop the Tumor / Analyte Type columns
# that were previously selected from smp.
$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND f1.experimental_strategy_and_data_subtypes = 'RNA-Seq|WXS'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value2 = $newSamplesQuery

# Clear the TsvExcel / WebExcel file-name references for the SamplesTab (row 3)
# and FilesTab (row 4) rows - these tabs no longer carry those values.
$ws.Range("D3:E4").ClearContents()

# Reflect the updated cursor/selection position left by the author.
$ws.Range("C3").Select() | Out-Null
